$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "www.google.com" + a/b/c/d header-like row
$ws.Range("A1").Formula = "=""www.google.com"""
$ws.Range("B1").Formula = "=""a"""
$ws.Range("C1").Formula = "=""b"""
$ws.Range("D1").Formula = "=""c"""
$ws.Range("E1").Formula = "=""d"""

# Row 2: same URL, then numeric-looking values that must remain TEXT, then two
# blank-but-present text cells. Using a quoted-string formula (instead of
# .Value) is what keeps "5"/"6" as text instead of being coerced to numbers.
$ws.Range("A2").Formula = "=""www.google.com"""
$ws.Range("B2").Formula = "=""5"""
$ws.Range("C2").Formula = "=""6"""
$ws.Range("D2").Formula = "="""""
$ws.Range("E2").Formula = "="""""

# Sheet reads left-to-right.
$excel.ActiveWindow.DisplayRightToLeft = $false
